# Japanese verb conjugation workbook update:
#  - add four new verb rows (送る / 上がる / 下る / 下がる groups) under the
#    existing 通う row
#  - fix a typo in the Nai Form of "通う" (D82) while filling in the new rows
#  - fill in the three previously-blank Volitional Form cells (F72:F74)
#  - add one more new verb row (集まる group)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the row-82 formatting (font / style / row height) down onto the new
# rows before typing the data into them, so every new cell picks up style
# "s=3" (Yu Gothic) just like the rest of the table.
$ws.Range("A82:F82").Copy()
$ws.Range("A83:F86").PasteSpecial(-4122)
$ws.Rows.Item(83).RowHeight = 18.75
$ws.Rows.Item(84).RowHeight = 18.75
$ws.Rows.Item(85).RowHeight = 18.75
$ws.Rows.Item(86).RowHeight = 18.75

# Dictionary Form column first for the four new verbs.
$ws.Range("A83").Value = "送る"
$ws.Range("A84").Value = "上がる"
$ws.Range("A85").Value = "下る"
$ws.Range("A86").Value = "下がる"

# 送る row (83): Te / Ta Form.
$ws.Range("B83").Value = "送って"
$ws.Range("C83").Value = "送った"

# Fix the typo'd Nai Form of 通う spotted while filling the Nai Form column.
$ws.Range("D82").Value = "通わない"

# 送る row (83): Nai / Masu / Volitional Form.
$ws.Range("D83").Value = "送らない"
$ws.Range("E83").Value = "送ります"
$ws.Range("F83").Value = "送ろう"

# 上がる row (84).
$ws.Range("B84").Value = "上がって"
$ws.Range("C84").Value = "上がった"
$ws.Range("D84").Value = "上がらない"
$ws.Range("E84").Value = "上がります"
$ws.Range("F84").Value = "上がろう"

# 下る row (85).
$ws.Range("B85").Value = "下って"
$ws.Range("C85").Value = "下った"
$ws.Range("D85").Value = "下らない"
$ws.Range("E85").Value = "下ります"
$ws.Range("F85").Value = "下ろう"

# 下がる row (86).
$ws.Range("B86").Value = "下がって"
$ws.Range("C86").Value = "下がった"
$ws.Range("D86").Value = "下がらない"
$ws.Range("E86").Value = "下がります"
$ws.Range("F86").Value = "下がろう"

# Back up to fill in the three Volitional Form cells that were left as "na"
# placeholders, bottom-to-top.
$ws.Range("F74").Value = "始まろう"
$ws.Range("F73").Value = "困ろう"
$ws.Range("F72").Value = "謝ろう"

# Finally, add the 集まる row (87) as a single new row.
$ws.Range("A82:F82").Copy()
$ws.Range("A87:F87").PasteSpecial(-4122)
$ws.Rows.Item(87).RowHeight = 18.75

$ws.Range("A87").Value = "集まる"
$ws.Range("B87").Value = "集まって"
$ws.Range("C87").Value = "集まった"
$ws.Range("D87").Value = "集まらない"
$ws.Range("E87").Value = "集まります"
$ws.Range("F87").Value = "集まろう"

# Move the selection to match where the author ended up.
[void]$ws.Range("F87").Select()
